$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$range = $ws.Range("B1:B1048576")
$fcs = $range.FormatConditions

while ($fcs.Count -gt 0) {
    $fcs.Item(1).Delete()
}

$r1 = $fcs.Add(8, 3, '="no comenzado"')
$r1.Interior.Color = 255
$r2 = $fcs.Add(8, 3, '="en proceso"')
$r2.Interior.Color = 65535
$r3 = $fcs.Add(8, 3, '="terminado"')
$r3.Interior.Color = 5296274
